# Update odds values on the "Jogos da Semana" sheet to reflect the latest
# FlashScore data refresh, per the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ind. Medellin vs Jaguares de Cordoba)
$ws.Range("G2").Value  = 1.55   # Odd_H_FT
$ws.Range("H2").Value  = 3.6    # Odd_D_FT
$ws.Range("I2").Value  = 6.5    # Odd_A_FT
$ws.Range("L2").Value  = 6.5    # Odd_A_HT
$ws.Range("Q2").Value  = 2.2    # Odd_Over25_FT
$ws.Range("R2").Value  = 1.65   # Odd_Under25_FT
$ws.Range("X2").Value  = 6.5    # Odd_CS_2-0
$ws.Range("Z2").Value  = 11     # Odd_CS_3-0
$ws.Range("AG2").Value = 13     # Odd_CS_0-1
$ws.Range("AJ2").Value = 67     # Odd_CS_0-3
$ws.Range("AL2").Value = 51     # Odd_CS_2-3
$ws.Range("AN2").Value = 3.4    # Odd_CS_1-0_HT
$ws.Range("AX2").Value = 34     # Odd_CS_0-2_HT

# Row 3 (Seattle Sounders vs Houston Dynamo)
$ws.Range("Q3").Value = 2.03    # Odd_Over25_FT
$ws.Range("R3").Value = 1.83    # Odd_Under25_FT
